$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Remove the duplicate "VIEJO RIVAS MAYRA ANABELLE" row (row 7); this shifts the
# trailing summary row (old row 8) up to become row 7.
$ws1.Rows.Item(7).Delete()

# Row 6's client becomes "VIEJO RIVAS MAYRA ANABELLE" (was "LATACELA ZUÑIGA JUAN FERNANDO")
$ws1.Range("B6").Value = "VIEJO RIVAS MAYRA ANABELLE"

# The summary row (now row 7) counted 6 data rows before ("0 de 6"); now there
# are only 5 data rows left, so it becomes "0 de 5".
$ws1.Range("C7:R7").Value = "0 de 5"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same row removal as sheet 1.
$ws2.Rows.Item(7).Delete()

# Row 6's client becomes "VIEJO RIVAS MAYRA ANABELLE" and its sales figures reset to 0.
$ws2.Range("B6").Value = "VIEJO RIVAS MAYRA ANABELLE"
$ws2.Range("C6").Value = 0
$ws2.Range("G6").Value = 0

# The totals row (now row 7) recalculates to all zeros since every data row is 0.
$ws2.Range("C7:G7").Value = 0

# Column C width narrows from 12 to 10 (stored OOXML width). The COM
# ColumnWidth property is offset from the stored width by 5/6, so we need to
# request 10 - 5/6 to land exactly on a stored width of 10.
$ws2.Columns.Item(3).ColumnWidth = 9.166666666666666
